$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# ---- Title shape: merge " " + "Sep 2018 DSUC Lab Practicals" into one run ----
$titleShape = $s.Shapes.Item(1)
$titleTr = $titleShape.TextFrame.TextRange
$titlePara1 = $titleTr.Paragraphs(1, 1)
$suffix = " Sep 2018 DSUC Lab Practicals"
$startPos = $titlePara1.Text.IndexOf($suffix) + 1
$titleSub = $titlePara1.Characters($startPos, $suffix.Length)
$titleSub.Text = $suffix

# ---- Content shape ----
$contentShape = $s.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

# Bump the numbered-list start value from 15 to 16
$para1 = $tr.Paragraphs(1, 1)
$para1.ParagraphFormat.Bullet.StartValue = 16

# Merge "doubly linked " + "list" into a single run "doubly linked list"
$dlText = "doubly linked list"
$dlPos = $para1.Text.IndexOf($dlText) + 1
$dl = $para1.Characters($dlPos, $dlText.Length)
$dl.Text = $dlText

# Remove the "Deletion at beginning / end / given position" bullet paragraphs
# (paragraph 5 three times in a row - each Delete() shifts later paragraphs up)
$tr.Paragraphs(5, 1).Delete()
$tr.Paragraphs(5, 1).Delete()
$tr.Paragraphs(5, 1).Delete()

# The "After each operation, display (traverse) the linked list." paragraph
# is now paragraph 5. Insert a blank paragraph before it (inherits its pPr).
$afterPara = $tr.Paragraphs(5, 1)
$afterPara.InsertBefore("`r") | Out-Null

# Split "After each operation..." text into two runs: "After " / "each operation..."
$afterPara = $tr.Paragraphs(6, 1)
$afterPrefix = $afterPara.Characters(1, 6)
$afterPrefix.Text = "After "
